# Auto-generated edit script to update column F (想去人数 / want-to-go count) values
# across all 4 worksheets, per the commit "Update gh-pages to output generated at 456a3b4"
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(3, 6).Value2 = 953   # F3: 950 -> 953
$ws.Cells.Item(4, 6).Value2 = 56   # F4: 53 -> 56
$ws.Cells.Item(7, 6).Value2 = 1172   # F7: 1171 -> 1172
$ws.Cells.Item(8, 6).Value2 = 938   # F8: 937 -> 938
$ws.Cells.Item(9, 6).Value2 = 33   # F9: 31 -> 33
$ws.Cells.Item(11, 6).Value2 = 1046   # F11: 1043 -> 1046
$ws.Cells.Item(12, 6).Value2 = 1492   # F12: 1489 -> 1492
$ws.Cells.Item(15, 6).Value2 = 1663   # F15: 1660 -> 1663
$ws.Cells.Item(17, 6).Value2 = 630   # F17: 628 -> 630
$ws.Cells.Item(21, 6).Value2 = 1091   # F21: 1090 -> 1091
$ws.Cells.Item(22, 6).Value2 = 1520   # F22: 1517 -> 1520
$ws.Cells.Item(23, 6).Value2 = 764   # F23: 763 -> 764
$ws.Cells.Item(24, 6).Value2 = 635   # F24: 634 -> 635
$ws.Cells.Item(26, 6).Value2 = 481   # F26: 480 -> 481
$ws.Cells.Item(28, 6).Value2 = 22   # F28: 17 -> 22
$ws.Cells.Item(29, 6).Value2 = 1026   # F29: 1024 -> 1026
$ws.Cells.Item(30, 6).Value2 = 1159   # F30: 1156 -> 1159
$ws.Cells.Item(32, 6).Value2 = 2444   # F32: 2439 -> 2444
$ws.Cells.Item(33, 6).Value2 = 282   # F33: 281 -> 282
$ws.Cells.Item(34, 6).Value2 = 1403   # F34: 1399 -> 1403
$ws.Cells.Item(36, 6).Value2 = 3   # F36: 2 -> 3
$ws.Cells.Item(38, 6).Value2 = 4012   # F38: 4008 -> 4012

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(4, 6).Value2 = 1040   # F4: 1039 -> 1040
$ws.Cells.Item(6, 6).Value2 = 191   # F6: 188 -> 191
$ws.Cells.Item(14, 6).Value2 = 4139   # F14: 4138 -> 4139
$ws.Cells.Item(26, 6).Value2 = 238   # F26: 237 -> 238

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(5, 6).Value2 = 1678   # F5: 1675 -> 1678
$ws.Cells.Item(7, 6).Value2 = 1026   # F7: 1024 -> 1026

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(3, 6).Value2 = 1678   # F3: 1675 -> 1678
$ws.Cells.Item(5, 6).Value2 = 1026   # F5: 1024 -> 1026
$ws.Cells.Item(7, 6).Value2 = 953   # F7: 950 -> 953
$ws.Cells.Item(8, 6).Value2 = 56   # F8: 53 -> 56
$ws.Cells.Item(10, 6).Value2 = 1172   # F10: 1171 -> 1172
$ws.Cells.Item(11, 6).Value2 = 938   # F11: 937 -> 938
$ws.Cells.Item(13, 6).Value2 = 33   # F13: 31 -> 33
$ws.Cells.Item(15, 6).Value2 = 191   # F15: 188 -> 191
$ws.Cells.Item(16, 6).Value2 = 191   # F16: 188 -> 191
$ws.Cells.Item(17, 6).Value2 = 1046   # F17: 1043 -> 1046
$ws.Cells.Item(19, 6).Value2 = 1492   # F19: 1489 -> 1492
$ws.Cells.Item(21, 6).Value2 = 1663   # F21: 1660 -> 1663
$ws.Cells.Item(23, 6).Value2 = 630   # F23: 628 -> 630
$ws.Cells.Item(26, 6).Value2 = 1520   # F26: 1517 -> 1520
$ws.Cells.Item(28, 6).Value2 = 764   # F28: 763 -> 764
$ws.Cells.Item(29, 6).Value2 = 635   # F29: 634 -> 635
$ws.Cells.Item(31, 6).Value2 = 481   # F31: 480 -> 481
$ws.Cells.Item(32, 6).Value2 = 22   # F32: 17 -> 22
$ws.Cells.Item(36, 6).Value2 = 1026   # F36: 1024 -> 1026
$ws.Cells.Item(40, 6).Value2 = 2445   # F40: 2439 -> 2445
$ws.Cells.Item(41, 6).Value2 = 238   # F41: 237 -> 238
$ws.Cells.Item(46, 6).Value2 = 1403   # F46: 1399 -> 1403
$ws.Cells.Item(47, 6).Value2 = 3   # F47: 2 -> 3
$ws.Cells.Item(49, 6).Value2 = 4012   # F49: 4008 -> 4012
